$d = $word.ActiveDocument

$replacements = @(
    @("152×5=", "593×8="),
    @("685×6=", "844×4="),
    @("141×4=", "257×8="),
    @("120×6=", "804×9="),
    @("742×6=", "490×4="),
    @("790×5=", "975×6="),
    @("980×7=", "185×2="),
    @("153×5=", "568×3="),
    @("572×4=", "854×3="),
    @("709×4=", "972×8="),
    @("441×3=", "596×5="),
    @("905×6=", "121×9="),
    @("317×9=", "308×6="),
    @("560×4=", "843×3="),
    @("247×6=", "592×6="),
    @("910×8=", "141×7="),
    @("481×6=", "928×7="),
    @("779×4=", "707×9="),
    @("968×8=", "225×9="),
    @("706×4=", "935×4="),
    @("808×5=", "943×3="),
    @("175×8=", "506×5="),
    @("256×6=", "263×2="),
    @("882×9=", "635×2="),
    @("138×9=", "787×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
